$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark (currently sits right after
#     the "runMCMC" run) before we create its replacement further up in the
#     document, so the name lookup below is unambiguous. ---
try {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
} catch {
    # no existing bookmark -- nothing to remove
}

# --- Step 2: rewrite the Heading-1 title paragraph. Text stays the same
#     ("Doing Bayesian Statistics Matlab Toolbox") but it now carries
#     en-US language formatting on the paragraph mark and on every run,
#     and the run split collapses to one run for the non-spellchecked
#     words plus a dedicated spell-check-wrapped run for "Matlab". ---
$p1 = $d.Paragraphs.Item(1)
$heading1Xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="berschrift1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Doing Bayesian Statistics </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Matlab</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Toolbox</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p1.Range.InsertXML($heading1Xml)

# --- Step 3: rewrite the 4th paragraph (the lone "- " bullet under
#     "Names") so it reads "- Matlab Toolbox for Bayesian Estimation
#     (MBE)" and carries the (relocated) "_GoBack" bookmark at its end. ---
$p4 = $d.Paragraphs.Item(4)
$namesLine2Xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Matlab</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Toolbox for Bayesian Estimation (MBE)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p4.Range.InsertXML($namesLine2Xml)

Write-Output "done"
